# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Rebuilds the mora-period table for the two workers: previously each
# worker had periods 1812, 1901, 1902 (x79434) and 1903 (x66195) spread
# across rows 16-23 in a mixed order. The new data lists each worker's
# four periods (1903, 1902, 1901, 1812) together, newest first, with the
# 66195 value following the 1903 period for both workers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Worker 1: HECTOR OSWALDO GARCIA MOLINA (CC 1051636273) -> rows 16-19
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 66195

$ws.Range("C17").Value = "1051636273"
$ws.Range("D17").Value = "HECTOR OSWALDO GARCIA MOLINA"
$ws.Range("E17").Value = "1902"

$ws.Range("E19").Value = "1812"
$ws.Range("C19").Value = "1051636273"
$ws.Range("D19").Value = "HECTOR OSWALDO GARCIA MOLINA"

# Worker 2: DENFRY MANUEL CASTRO BAENA (CC 1193115546) -> rows 20-23
$ws.Range("E20").Value = "1903"
$ws.Range("F20").Value = 66195

$ws.Range("C21").Value = "1193115546"
$ws.Range("D21").Value = "DENFRY MANUEL CASTRO BAENA"

$ws.Range("C22").Value = "1193115546"
$ws.Range("D22").Value = "DENFRY MANUEL CASTRO BAENA"
$ws.Range("E22").Value = "1901"
$ws.Range("F22").Value = 79434

$ws.Range("E23").Value = "1812"
$ws.Range("F23").Value = 79434
